# Invoices.xlsx edit:
#  - "Invoice Line Items" sheet gets a new leading "ID" column (1..19),
#    shifting the existing Invoice ID / Product ID / # of Items columns
#    one column to the right (B,C,D).
#  - Selection on "Invoice Line Items" moves to A21 (just past the data).
#  - "Invoices" sheet becomes the active tab/sheet, with its selection
#    moved to B6.

$wb = $excel.ActiveWorkbook

# --- "Invoice Line Items": insert an ID column in front of the data ---
$ws3 = $wb.Worksheets.Item("Invoice Line Items")

for ($r = 1; $r -le 20; $r++) {
    $oldA = $ws3.Cells.Item($r, 1).Value2
    $oldB = $ws3.Cells.Item($r, 2).Value2
    $oldC = $ws3.Cells.Item($r, 3).Value2

    $ws3.Cells.Item($r, 4).Value = $oldC
    $ws3.Cells.Item($r, 3).Value = $oldB
    $ws3.Cells.Item($r, 2).Value = $oldA
}

$ws3.Cells.Item(1, 1).Value = "ID"
for ($r = 2; $r -le 20; $r++) {
    $ws3.Cells.Item($r, 1).Value = $r - 1
}

$ws3.Range("A21").Select() | Out-Null

# --- "Invoices": make it the active sheet/tab, move selection to B6 ---
$ws2 = $wb.Worksheets.Item("Invoices")
$ws2.Activate() | Out-Null
$ws2.Range("B6").Select() | Out-Null
